$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Top table (rows 2-7): change the "x" (periode) column B2:B7. Formulas in
# D (xy) and E (x2) as well as the summary row 8 and the A12/A15/A18 results
# recompute automatically.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 5

# ---------------------------------------------------------------------------
# Bottom table (rows 22-27): update "penjualan" (B), "peramalan" (C) and the
# "|at-ft|" (E) columns. D (at-ft) is a formula and recalculates on its own.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 4
$ws.Range("E22").Value = 1

$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = 2

$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 4
$ws.Range("E24").Value = 0

$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 1

$ws.Range("B26").Value = 3
$ws.Range("C26").Value = 4
$ws.Range("E26").Value = 1

$ws.Range("B27").Value = 5
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 1

# ---------------------------------------------------------------------------
# New explanatory notes in column G (added in the order the original author
# typed them so the shared-string table comes out in the same sequence).
# ---------------------------------------------------------------------------
$ws.Range("G31").Value = "ket : total |at-ft| dibagi lama periode yang terlibat"
$ws.Range("G32").Value = "peramalan didapat dari hasil peramalan yang didapat dari proses sebelumnya"
$ws.Range("G30").Value = "at : penjualan, ft : peramalan"
$ws.Range("G22").Value = "<-- inget ganti bagian ini"

$ws.Range("G31").WrapText = $true
$ws.Range("G32").WrapText = $true
$ws.Rows.Item(31).RowHeight = 30
$ws.Rows.Item(32).RowHeight = 60

# ---------------------------------------------------------------------------
# Highlight the revised |at-ft| column with a yellow fill + thin box border.
# ---------------------------------------------------------------------------
$rng = $ws.Range("E22:E27")
$rng.Borders.LineStyle = 1
$rng.Interior.Color = 65535

# ---------------------------------------------------------------------------
# Move the selection to match the author's last position.
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("G27"))

Write-Output "done"
